$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Info")

# 1) Rename player "topdosl33ts" -> "CarvaPower." (row 25, column A)
$ws.Cells.Item(25, 1).Value = "CarvaPower."

# 2) Swap the two whole rows that held "Marcelo" (row 37) and "supercell"
#    (row 38): the name (A) and every stat column (B-I) trade places.
$marcelo = @($ws.Cells.Item(37, 1).Value2, $ws.Cells.Item(37, 2).Value2, $ws.Cells.Item(37, 3).Value2, $ws.Cells.Item(37, 4).Value2, $ws.Cells.Item(37, 5).Value2, $ws.Cells.Item(37, 6).Value2, $ws.Cells.Item(37, 7).Value2, $ws.Cells.Item(37, 8).Value2, $ws.Cells.Item(37, 9).Value2)
$supercell = @($ws.Cells.Item(38, 1).Value2, $ws.Cells.Item(38, 2).Value2, $ws.Cells.Item(38, 3).Value2, $ws.Cells.Item(38, 4).Value2, $ws.Cells.Item(38, 5).Value2, $ws.Cells.Item(38, 6).Value2, $ws.Cells.Item(38, 7).Value2, $ws.Cells.Item(38, 8).Value2, $ws.Cells.Item(38, 9).Value2)

for ($i = 0; $i -lt 9; $i++) {
    $ws.Cells.Item(37, 1 + $i).Value = $supercell[$i]
}
for ($i = 0; $i -lt 9; $i++) {
    $ws.Cells.Item(38, 1 + $i).Value = $marcelo[$i]
}

# The WARSCORE (J) column keeps its own arithmetic: 62 -> 61 for row 37,
# while row 38's J value (61) is unchanged, so set both explicitly.
$ws.Cells.Item(37, 10).Value = 61.0
$ws.Cells.Item(38, 10).Value = 61.0

# 3) WARSCORE (column J) bumps of +/-1 on several other rows.
$ws.Cells.Item(14, 10).Value = 202.0
$ws.Cells.Item(19, 10).Value = 180.0
$ws.Cells.Item(25, 10).Value = 131.0
$ws.Cells.Item(26, 10).Value = 130.0
$ws.Cells.Item(34, 10).Value = 76.0
$ws.Cells.Item(40, 10).Value = 44.0
$ws.Cells.Item(51, 10).Value = -71.0
